$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price column cells whose new values look numeric,
# so Excel keeps them as text (matching the original inlineStr string cells)
# instead of silently converting to a Number and losing formatting (e.g. trailing zeros).
$ws.Range("D2").Value = "24.451.15"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.654.66"
$ws.Range("E3").Value = "  -2.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.51"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9981"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3625"
$ws.Range("E7").Value = "  -2.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.46"
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3260"
$ws.Range("E9").Value = "  -4.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.123"
$ws.Range("E10").Value = "  -4.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06959"
$ws.Range("E11").Value = "  -6.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.917"
$ws.Range("E13").Value = "  -4.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.31"
$ws.Range("E14").Value = "  -6.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.598"
$ws.Range("E15").Value = "  -4.10%  "
$ws.Range("D16").Value = "1.652.93"
$ws.Range("E16").Value = "  -2.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001043"
$ws.Range("E17").Value = "  -6.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06525"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9985"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "76.37"
$ws.Range("E20").Value = "  -7.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.902"
$ws.Range("E21").Value = "  -6.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.70"
$ws.Range("E22").Value = "  -7.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.52"
$ws.Range("E23").Value = "  -4.96%  "
$ws.Range("D24").Value = "24.453.07"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.452"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.299"
$ws.Range("E26").Value = "  -16.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.82"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.47"
$ws.Range("E28").Value = "  -7.88%  "
$ws.Range("D29").Value = "1.840.80"
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.34"
$ws.Range("E30").Value = "  -4.61%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.187"
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.589"
$ws.Range("E33").Value = "  -15.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08338"
$ws.Range("E34").Value = "  -4.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.683"
$ws.Range("E35").Value = "  -4.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.34"
$ws.Range("E36").Value = "  -8.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.200"
$ws.Range("E37").Value = "  -4.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06049"
$ws.Range("E38").Value = "  -6.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02190"
$ws.Range("E39").Value = "  -6.92%  "
$ws.Range("E40").Value = "  -5.21%  "
$ws.Range("E41").Value = "  -5.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.154"
$ws.Range("E42").Value = "  -8.06%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5883"
$ws.Range("E44").Value = "  -7.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.732"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.65"
$ws.Range("E46").Value = "  -8.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5590"
$ws.Range("E47").Value = "  -7.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.05"
$ws.Range("E48").Value = "  -4.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.935"
$ws.Range("E49").Value = "  -7.44%  "
$ws.Range("E50").Value = "  -4.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.00"
$ws.Range("E51").Value = "  -5.70%  "
